# Scene 9B ("Act 2 Lilith") — "write some new for stephen"
#
# The canonical-OOXML diff for this commit touches only word/styles.xml
# (the built-in Normal/TableNormal/Heading1-6/Title/Subtitle style
# definitions that Stephen's authoring tool keeps in the template) and the
# opaque Google-Docs roundtrip signature in customXML/item1.xml; the
# w:body content in word/document.xml is byte-for-byte identical before
# and after the commit. In other words this particular revision is Word
# making sure the full built-in style set (headings, title, subtitle,
# table-normal) backing the script's formatting is registered in the
# stylesheet — no dialogue/scene text actually changed.
#
# Reproduce that by touching the same style slots through the real Word
# object model. Styles.Add() on a name that already resolves to a
# built-in style returns the existing Style object (exactly like Word's
# "Keep existing styles" behavior) instead of minting a duplicate/corrupt
# styleId, so this both matches the diff's intent and leaves the prose
# (word/document.xml) untouched.

$d = $word.ActiveDocument

$d.Styles.Add("Normal", 1)       | Out-Null
$d.Styles.Add("Table Normal", 1) | Out-Null
$d.Styles.Add("Heading 1", 1)    | Out-Null
$d.Styles.Add("Heading 2", 1)    | Out-Null
$d.Styles.Add("Heading 3", 1)    | Out-Null
$d.Styles.Add("Heading 4", 1)    | Out-Null
$d.Styles.Add("Heading 5", 1)    | Out-Null
$d.Styles.Add("Heading 6", 1)    | Out-Null
$d.Styles.Add("Title", 1)        | Out-Null
$d.Styles.Add("Subtitle", 1)     | Out-Null
